$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting the existing rows 134-174 down to 135-175.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data record.
$ws.Range("A134").Value = 10
$ws.Range("B134").Value = "Vega Modelo de Temuco"
$ws.Range("C134").Value = "La Araucanía"
$ws.Range("D134").Value = 45093
$ws.Range("E134").Value = 9
$ws.Range("F134").Value = 100112035
$ws.Range("G134").Value = "Bruselas (repollito)"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 35
$ws.Range("K134").Value = 28000
$ws.Range("L134").Value = 28000
$ws.Range("M134").Value = 28000
$ws.Range("N134").Value = "$/malla 15 kilos"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 1867
$ws.Range("Q134").Value = 15
$ws.Range("R134").Value = "Hortaliza"
